# Update the build timestamp embedded in the workbook's version strings.
#
# The workbook was re-built; the only change is that every occurrence of the
# old "built on" timestamp is replaced with the new one, across the "About"
# sheet (A2, A6) and the "Boundaries and methane sources" sheet (build_version
# column, rows 2-10).

$oldTimestamp = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newTimestamp = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- About sheet -----------------------------------------------------------
# A2: "Version: mines - January 30 (built on ...)"
$a2 = $aboutSheet.Range("A2").Value()
$aboutSheet.Range("A2").Value = $a2.Replace($oldTimestamp, $newTimestamp)

# A6: "Recommended Citation: ... version '...'. ..."
$a6 = $aboutSheet.Range("A6").Value()
$aboutSheet.Range("A6").Value = $a6.Replace($oldTimestamp, $newTimestamp)

# --- Boundaries and methane sources sheet ----------------------------------
# Column S ("build_version") holds the same version string for every data
# row (rows 2 through 10 in this workbook).
$usedRange = $dataSheet.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)   # column S = 19
    $val = $cell.Value()
    if ($val -ne $null -and $val.Contains($oldTimestamp)) {
        $cell.Value = $val.Replace($oldTimestamp, $newTimestamp)
    }
}
